$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 24 March 2022 FedEx shipment test-run results: refresh the
# ShipmentTracking (P), ActualRate (Q) and Result (R) columns for
# rows 2-25 with the latest captured values.

# Force columns P and Q to a text format before writing so that the
# numeric-looking tracking numbers and "$"-prefixed rate strings are
# stored as literal text, then restore the default "Normal" style so
# no extra cell formatting is left behind.
$ws.Range("P2:Q25").NumberFormat = "@"

$ws.Range("P2").Value = "320018113690"
$ws.Range("Q2").Value = "`$19.04"
$ws.Range("P3").Value = "320018113704"
$ws.Range("Q3").Value = "`$27.50"
$ws.Range("P4").Value = "320018113737"
$ws.Range("Q4").Value = "`$31.73"
$ws.Range("P5").Value = "320018113759"
$ws.Range("Q5").Value = "`$43.36"
$ws.Range("P6").Value = "320018113792"
$ws.Range("Q6").Value = "`$56.05"
$ws.Range("P7").Value = "320018113818"
$ws.Range("Q7").Value = "`$231.08"
$ws.Range("P8").Value = "320018113840"
$ws.Range("Q8").Value = "`$19.04"
$ws.Range("P9").Value = "320018113862"
$ws.Range("Q9").Value = "`$23.27"
$ws.Range("P10").Value = "320018113895"
$ws.Range("Q10").Value = "`$27.50"
$ws.Range("P11").Value = "320018113910"
$ws.Range("Q11").Value = "`$40.19"
$ws.Range("P12").Value = "320018113954"
$ws.Range("Q12").Value = "`$52.88"
$ws.Range("P13").Value = "320018113976"
$ws.Range("Q13").Value = "`$14.81"
$ws.Range("P14").Value = "320018114001"
$ws.Range("Q14").Value = "`$17.98"
$ws.Range("P15").Value = "320018114023"
$ws.Range("Q15").Value = "`$21.15"
$ws.Range("P16").Value = "320018114056"
$ws.Range("Q16").Value = "`$31.73"
$ws.Range("P17").Value = "320018114078"
$ws.Range("Q17").Value = "`$42.30"
$ws.Range("P18").Value = "320018114115"
$ws.Range("Q18").Value = "`$53.93"
$ws.Range("P19").Value = "320018114137"
$ws.Range("Q19").Value = "`$53.93"
$ws.Range("P20").Value = "320018114160"
$ws.Range("Q20").Value = "`$85.66"
$ws.Range("P21").Value = "320018114181"
$ws.Range("Q21").Value = "`$111.04"
$ws.Range("P22").Value = "320018114218"
$ws.Range("Q22").Value = "`$248.35"
$ws.Range("P23").Value = "320018115144"
$ws.Range("Q23").Value = "`$492.15"
$ws.Range("P24").Value = "320018115155"
$ws.Range("Q24").Value = "`$354.26"
$ws.Range("P25").Value = "320018115166"
$ws.Range("Q25").Value = "`$132.19"

$ws.Range("P2:Q25").Style = "Normal"

# Rows 2-22 passed this run, so their Result cell is cleared; rows
# 23-25 are the new failing shipments.
$ws.Range("R2").ClearContents() | Out-Null
$ws.Range("R3").ClearContents() | Out-Null
$ws.Range("R4").ClearContents() | Out-Null
$ws.Range("R5").ClearContents() | Out-Null
$ws.Range("R6").ClearContents() | Out-Null
$ws.Range("R7").ClearContents() | Out-Null
$ws.Range("R8").ClearContents() | Out-Null
$ws.Range("R9").ClearContents() | Out-Null
$ws.Range("R10").ClearContents() | Out-Null
$ws.Range("R11").ClearContents() | Out-Null
$ws.Range("R12").ClearContents() | Out-Null
$ws.Range("R13").ClearContents() | Out-Null
$ws.Range("R14").ClearContents() | Out-Null
$ws.Range("R15").ClearContents() | Out-Null
$ws.Range("R16").ClearContents() | Out-Null
$ws.Range("R17").ClearContents() | Out-Null
$ws.Range("R18").ClearContents() | Out-Null
$ws.Range("R19").ClearContents() | Out-Null
$ws.Range("R20").ClearContents() | Out-Null
$ws.Range("R21").ClearContents() | Out-Null
$ws.Range("R22").ClearContents() | Out-Null

$ws.Range("R23:R25").NumberFormat = "@"
$ws.Range("R23").Value = "FAIL"
$ws.Range("R24").Value = "FAIL"
$ws.Range("R25").Value = "FAIL"
$ws.Range("R23:R25").Style = "Normal"

# Match the sheet view selection recorded after the refresh.
$ws.Range("R2:R25").Select() | Out-Null
